$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed crypto price / 1h-volume data (GitHub Actions scrape).
# Every target cell in the sheet is stored as a string (t="s" / t="inlineStr"),
# so values are written with a leading apostrophe to force text storage even when
# the text looks numeric (e.g. "1.00"), then the style is reset to Normal so Excel
# does not leave a quote-prefix / number-format style attached to the cell (matching
# the original, unstyled data cells).

$ws.Range("D2").Value = "'26.626.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.596.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.41%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.18%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.15%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.0618"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.37%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.06%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.14%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.820.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.584.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.65%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.33%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.619.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.04%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0737"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.76%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'208.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +5.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.80%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'145.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.81%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.50%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.40%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.07%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.275.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.28%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.619"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -7.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.10%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.85%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.55%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.968"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +16.37%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.61%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.94%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'64.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.67%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.732.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.17%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'90.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.57%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.18%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'USDD"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.12%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'7.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.82%  "
$ws.Range("E51").Style = "Normal"
